$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 486, shifting the existing data (rows 486-598)
# down to rows 490-602. This also carries formatting (e.g. the date style on
# column D) down onto the newly inserted blank rows.
$ws.Rows("486:489").Insert()

# Fill the 4 newly inserted rows with this week's price data.
# Columns A,B,C,E,F,G,H,I,J,K,Q,R,T are constant for this market/product
# subset, matching the rest of the sheet.

# Row 486: Femacal de La Calera / Frutilla - Especial
$ws.Range("A486").Value = 3
$ws.Range("B486").Value = "Femacal de La Calera"
$ws.Range("C486").Value = "Coquimbo"
$ws.Range("D486").Value = 45211
$ws.Range("E486").Value = 5
$ws.Range("F486").Value = "Fruta"
$ws.Range("G486").Value = 100101
$ws.Range("H486").Value = "Berries"
$ws.Range("I486").Value = 100112025
$ws.Range("J486").Value = "Frutilla"
$ws.Range("K486").Value = "Sin especificar"
$ws.Range("L486").Value = "Especial"
$ws.Range("M486").Value = 138
$ws.Range("N486").Value = 11000
$ws.Range("O486").Value = 12000
$ws.Range("P486").Value = 11565
$ws.Range("Q486").Value = "$/bandeja 7 kilos"
$ws.Range("R486").Value = "Provincia de Melipilla"
$ws.Range("S486").Value = 1652
$ws.Range("T486").Value = 7

# Row 487: Primera
$ws.Range("A487").Value = 3
$ws.Range("B487").Value = "Femacal de La Calera"
$ws.Range("C487").Value = "Coquimbo"
$ws.Range("D487").Value = 45211
$ws.Range("E487").Value = 5
$ws.Range("F487").Value = "Fruta"
$ws.Range("G487").Value = 100101
$ws.Range("H487").Value = "Berries"
$ws.Range("I487").Value = 100112025
$ws.Range("J487").Value = "Frutilla"
$ws.Range("K487").Value = "Sin especificar"
$ws.Range("L487").Value = "Primera"
$ws.Range("M487").Value = 125
$ws.Range("N487").Value = 8000
$ws.Range("O487").Value = 9000
$ws.Range("P487").Value = 8600
$ws.Range("Q487").Value = "$/bandeja 7 kilos"
$ws.Range("R487").Value = "Provincia de Melipilla"
$ws.Range("S487").Value = 1229
$ws.Range("T487").Value = 7

# Row 488: Segunda
$ws.Range("A488").Value = 3
$ws.Range("B488").Value = "Femacal de La Calera"
$ws.Range("C488").Value = "Coquimbo"
$ws.Range("D488").Value = 45211
$ws.Range("E488").Value = 5
$ws.Range("F488").Value = "Fruta"
$ws.Range("G488").Value = 100101
$ws.Range("H488").Value = "Berries"
$ws.Range("I488").Value = 100112025
$ws.Range("J488").Value = "Frutilla"
$ws.Range("K488").Value = "Sin especificar"
$ws.Range("L488").Value = "Segunda"
$ws.Range("M488").Value = 56
$ws.Range("N488").Value = 6000
$ws.Range("O488").Value = 6000
$ws.Range("P488").Value = 6000
$ws.Range("Q488").Value = "$/bandeja 7 kilos"
$ws.Range("R488").Value = "Provincia de Melipilla"
$ws.Range("S488").Value = 857
$ws.Range("T488").Value = 7

# Row 489: Tercera
$ws.Range("A489").Value = 3
$ws.Range("B489").Value = "Femacal de La Calera"
$ws.Range("C489").Value = "Coquimbo"
$ws.Range("D489").Value = 45211
$ws.Range("E489").Value = 5
$ws.Range("F489").Value = "Fruta"
$ws.Range("G489").Value = 100101
$ws.Range("H489").Value = "Berries"
$ws.Range("I489").Value = 100112025
$ws.Range("J489").Value = "Frutilla"
$ws.Range("K489").Value = "Sin especificar"
$ws.Range("L489").Value = "Tercera"
$ws.Range("M489").Value = 56
$ws.Range("N489").Value = 4000
$ws.Range("O489").Value = 4000
$ws.Range("P489").Value = 4000
$ws.Range("Q489").Value = "$/bandeja 7 kilos"
$ws.Range("R489").Value = "Provincia de Melipilla"
$ws.Range("S489").Value = 571
$ws.Range("T489").Value = 7
